{"js": "// Update the 24 \"dividend\u00f7divisor=\" cell values in the practice-problems\n// table to the new set of problems from the commit's target revision.\n// Each old value is unique in the document, so an exact (case-sensitive,\n// non-wildcard) search-and-replace on the body text is safe and\n// unambiguous.\nconst mapping = [\n  [\"153\u00f76=\", \"204\u00f79=\"],\n  [\"482\u00f79=\", \"593\u00f79=\"],\n  [\"536\u00f77=\", \"606\u00f76=\"],\n  [\"659\u00f76=\", \"601\u00f73=\"],\n  [\"913\u00f74=\", \"345\u00f72=\"],\n  [\"842\u00f72=\", \"742\u00f74=\"],\n  [\"660\u00f76=\", \"564\u00f73=\"],\n  [\"836\u00f74=\", \"870\u00f75=\"],\n  [\"229\u00f77=\", \"596\u00f76=\"],\n  [\"122\u00f79=\", \"811\u00f72=\"],\n  [\"304\u00f73=\", \"540\u00f75=\"],\n  [\"631\u00f72=\", \"947\u00f79=\"],\n  [\"479\u00f72=\", \"483\u00f76=\"],\n  [\"428\u00f72=\", \"802\u00f74=\"],\n  [\"976\u00f78=\", \"371\u00f79=\"],\n  [\"748\u00f72=\", \"878\u00f73=\"],\n  [\"579\u00f75=\", \"177\u00f78=\"],\n  [\"838\u00f74=\", \"829\u00f73=\"],\n  [\"735\u00f79=\", \"589\u00f75=\"],\n  [\"713\u00f74=\", \"695\u00f74=\"],\n  [\"393\u00f77=\", \"185\u00f74=\"],\n  [\"281\u00f76=\", \"450\u00f74=\"],\n  [\"860\u00f76=\", \"454\u00f75=\"],\n  [\"407\u00f75=\", \"613\u00f75=\"],\n  [\"429\u00f77=\", \"873\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of mapping) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n\n", "ps1": "# Update the 24 \"dividend/divisor=\" cell values in the practice-problems\n# table to the new set of problems from the commit's target revision.\n# Each old value is unique in the document, so Find/Replace on the whole\n# document body (Content) is safe and unambiguous - one hit per pair.\n\n$d = $word.ActiveDocument\n\n$mapping = @(\n    @{ Old = \"153\u00f76=\"; New = \"204\u00f79=\" },\n    @{ Old = \"482\u00f79=\"; New = \"593\u00f79=\" },\n    @{ Old = \"536\u00f77=\"; New = \"606\u00f76=\" },\n    @{ Old = \"659\u00f76=\"; New = \"601\u00f73=\" },\n    @{ Old = \"913\u00f74=\"; New = \"345\u00f72=\" },\n    @{ Old = \"842\u00f72=\"; New = \"742\u00f74=\" },\n    @{ Old = \"660\u00f76=\"; New = \"564\u00f73=\" },\n    @{ Old = \"836\u00f74=\"; New = \"870\u00f75=\" },\n    @{ Old = \"229\u00f77=\"; New = \"596\u00f76=\" },\n    @{ Old = \"122\u00f79=\"; New = \"811\u00f72=\" },\n    @{ Old = \"304\u00f73=\"; New = \"540\u00f75=\" },\n    @{ Old = \"631\u00f72=\"; New = \"947\u00f79=\" },\n    @{ Old = \"479\u00f72=\"; New = \"483\u00f76=\" },\n    @{ Old = \"428\u00f72=\"; New = \"802\u00f74=\" },\n    @{ Old = \"976\u00f78=\"; New = \"371\u00f79=\" },\n    @{ Old = \"748\u00f72=\"; New = \"878\u00f73=\" },\n    @{ Old = \"579\u00f75=\"; New = \"177\u00f78=\" },\n    @{ Old = \"838\u00f74=\"; New = \"829\u00f73=\" },\n    @{ Old = \"735\u00f79=\"; New = \"589\u00f75=\" },\n    @{ Old = \"713\u00f74=\"; New = \"695\u00f74=\" },\n    @{ Old = \"393\u00f77=\"; New = \"185\u00f74=\" },\n    @{ Old = \"281\u00f76=\"; New = \"450\u00f74=\" },\n    @{ Old = \"860\u00f76=\"; New = \"454\u00f75=\" },\n    @{ Old = \"407\u00f75=\"; New = \"613\u00f75=\" },\n    @{ Old = \"429\u00f77=\"; New = \"873\u00f73=\" }\n)\n\nforeach ($pair in $mapping) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $found) {\n        throw \"Text not found: $($pair.Old)\"\n    }\n}\n\n"}
